$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- Substantive content edits -------------------------------------------

# 1. Intro paragraph: bump discretewq package version and citation year.
Replace-Text `
    "The data integration code was packaged into the R package discretewq: https://github.com/sbashevkin/discretewq (Bashevkin 2021)." `
    "The data integration code was packaged into the R package discretewq v1.1.0: https://github.com/sbashevkin/discretewq (Bashevkin 2022)."

# 2. Data integration methods paragraph: same version/citation bump.
Replace-Text `
    "All data integration code can be found in the discretewq R package (https://github.com/sbashevkin/discretewq; Bashevkin 2021)." `
    "All data integration code can be found in the discretewq R package v1.1.0 (https://github.com/sbashevkin/discretewq; Bashevkin 2022)."

# 3. Literature cited entry for Bashevkin.
Replace-Text `
    "Bashevkin, S. M. 2021. discretewq: An Integrated Dataset of Discrete Water Quality in the San Francisco Estuary. Zenodo. doi:10.5281/zenodo.4631924" `
    "Bashevkin, S. M. 2022. discretewq: An Integrated Dataset of Discrete Water Quality in the San Francisco Estuary v1.1.0. Zenodo. doi:10.5281/zenodo.5834821"

# --- Cosmetic clean-up: collapse leftover spell/grammar-check run splits --
# (re-assert identical text across each split so Word recombines the runs
# and drops the now-stale w:proofErr markers, matching the authoritative
# re-save of the document)

Replace-Text `
    "Outside this short time period, this was measured as a" `
    "Outside this short time period, this was measured as a"

Replace-Text `
    "collected and analyzed similar to EMP." `
    "collected and analyzed similar to EMP."

Replace-Text `
    "From each dataset, we selected columns corresponding to the water quality variables of interest as well as important accessory information (date, time, station,  latitude, longitude, depth, tide, and any notes). We then renamed variables for consistency and converted all variables to consistent units. Salinity was calculated from specific conductivity using the ec2pss function from the wql R package (Jassby et al. 2017). This function uses the Practical Salinity Scale 1978 for salinities between 2 and 42 (Fofonoff and Millard Jr 1983) and the extension of the Practical Salinity Scale (Hill et al. 1986) for salinities below 2. Conductivity data were also retained in the integrated dataset. In most cases, latitude and longitude coordinates of the fixed sampling stations were retained. When these coordinates were not available (e.g. for non-fixed stations), we retained any coordinates that were recorded during the field sampling. To remove duplicate values from the dataset, only one set of values was retained for each recorded date, time, and location." `
    "From each dataset, we selected columns corresponding to the water quality variables of interest as well as important accessory information (date, time, station,  latitude, longitude, depth, tide, and any notes). We then renamed variables for consistency and converted all variables to consistent units. Salinity was calculated from specific conductivity using the ec2pss function from the wql R package (Jassby et al. 2017). This function uses the Practical Salinity Scale 1978 for salinities between 2 and 42 (Fofonoff and Millard Jr 1983) and the extension of the Practical Salinity Scale (Hill et al. 1986) for salinities below 2. Conductivity data were also retained in the integrated dataset. In most cases, latitude and longitude coordinates of the fixed sampling stations were retained. When these coordinates were not available (e.g. for non-fixed stations), we retained any coordinates that were recorded during the field sampling. To remove duplicate values from the dataset, only one set of values was retained for each recorded date, time, and location."

Replace-Text `
    "Fofonoff, N. P., and R. C. Millard Jr. 1983." `
    "Fofonoff, N. P., and R. C. Millard Jr. 1983."

Replace-Text `
    "Hill, K., T. Dauphinee, and D. Woods. 1986." `
    "Hill, K., T. Dauphinee, and D. Woods. 1986."

Replace-Text `
    "Jassby, A. D., J. E. Cloern, and J. Stachelek. 2017. wql: Exploring Water Quality Monitoring Data." `
    "Jassby, A. D., J. E. Cloern, and J. Stachelek. 2017. wql: Exploring Water Quality Monitoring Data."

Replace-Text `
    "Cloern, J. E., and T. S. Schraga. 2016. USGS Measurements of Water Quality in San Francisco Bay (CA), 1969-2015 (ver. 3.0 June 2017). U. S. Geological Survey data release. doi:https://doi.org/10.5066/F7TQ5ZPR" `
    "Cloern, J. E., and T. S. Schraga. 2016. USGS Measurements of Water Quality in San Francisco Bay (CA), 1969-2015 (ver. 3.0 June 2017). U. S. Geological Survey data release. doi:https://doi.org/10.5066/F7TQ5ZPR"

Replace-Text `
    "Interagency Ecological Program (IEP), L. Damon, T. Tempel, and A. Chorazyczewski. 2020a. Interagency Ecological Program San Francisco Estuary 20mm Survey 1995 - 2020. Environmental Data Initiative. doi:" `
    "Interagency Ecological Program (IEP), L. Damon, T. Tempel, and A. Chorazyczewski. 2020a. Interagency Ecological Program San Francisco Estuary 20mm Survey 1995 - 2020. Environmental Data Initiative. doi:"

Replace-Text `
    "Interagency Ecological Program (IEP), L. Damon, T. Tempel, and A. Chorazyczewski. 2020b. Interagency Ecological Program San Francisco Estuary Spring Kodiak Trawl Survey 2002 - 2020. Environmental Data Initiative. doi:" `
    "Interagency Ecological Program (IEP), L. Damon, T. Tempel, and A. Chorazyczewski. 2020b. Interagency Ecological Program San Francisco Estuary Spring Kodiak Trawl Survey 2002 - 2020. Environmental Data Initiative. doi:"

Replace-Text `
    "Schraga, T. S., E. S. Nejad, C. A. Martin, and J. E. Cloern. 2018. USGS measurements of water quality in San Francisco Bay (CA), beginning in 2016 (ver. 3.0, March 2020). U. S. Geological Survey data release. doi:https://doi.org/10.5066/F7D21WGF" `
    "Schraga, T. S., E. S. Nejad, C. A. Martin, and J. E. Cloern. 2018. USGS measurements of water quality in San Francisco Bay (CA), beginning in 2016 (ver. 3.0, March 2020). U. S. Geological Survey data release. doi:https://doi.org/10.5066/F7D21WGF"

Replace-Text `
    "United States Fish And Wildlife Service, C. Johnston, S. Durkacz, and others. 2020. Interagency Ecological Program and US Fish and Wildlife Service: San Francisco Estuary Enhanced Delta Smelt Monitoring Program data, 2016-2020. Environmental Data Initiative. doi:" `
    "United States Fish And Wildlife Service, C. Johnston, S. Durkacz, and others. 2020. Interagency Ecological Program and US Fish and Wildlife Service: San Francisco Estuary Enhanced Delta Smelt Monitoring Program data, 2016-2020. Environmental Data Initiative. doi:"

Replace-Text `
    "USBR, R. Dahlgren, L. Loken, and E. Van Nieuwenhuyse. 2020. Monthly vertical profiles of water quality in the Sacramento Deep Water Ship Channel 2012-2019." `
    "USBR, R. Dahlgren, L. Loken, and E. Van Nieuwenhuyse. 2020. Monthly vertical profiles of water quality in the Sacramento Deep Water Ship Channel 2012-2019."
